# Refresh the cryptocurrency price / 1h-volume figures (GitHub Actions scrape).
# Source diff changes only cell VALUES (cols B-E), rows 2-51; row 1 header and
# column A rank index are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new Price values are plain digit-dot-digit strings (e.g. "226.04",
# "1.00") that Excel would otherwise auto-coerce to numbers on assignment. Force
# those specific cells to Text format first so they stay strings, same as the rest
# of column D / the original file.
$textCells = @("D5", "D6", "D8", "D9", "D11", "D13", "D18", "D19", "D21", "D23", "D25", "D26", "D27", "D29", "D30", "D38", "D40", "D41", "D43", "D44", "D45", "D46", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Cell -> new value, in the same order as the diff.
$updates = New-Object System.Collections.Specialized.OrderedDictionary
$updates["D2"] = '34.611.30'
$updates["E2"] = '  +0.44%  '
$updates["D3"] = '1.811.70'
$updates["E3"] = '  +0.49%  '
$updates["E4"] = '  -0.23%  '
$updates["D5"] = '226.04'
$updates["E5"] = '  -1.07%  '
$updates["D6"] = '0.601'
$updates["E6"] = '  +3.65%  '
$updates["E7"] = '  -0.19%  '
$updates["D8"] = '36.57'
$updates["E8"] = '  +5.61%  '
$updates["D9"] = '0.293'
$updates["E9"] = '  -2.37%  '
$updates["D11"] = '0.0967'
$updates["E11"] = '  +1.54%  '
$updates["D12"] = '2.072.85'
$updates["E12"] = '  +0.49%  '
$updates["D13"] = '11.35'
$updates["E13"] = '  +1.15%  '
$updates["D14"] = '1.808.16'
$updates["E14"] = '  +0.29%  '
$updates["E15"] = '  -1.55%  '
$updates["D16"] = '34.553.63'
$updates["E17"] = '  +2.12%  '
$updates["D18"] = '68.75'
$updates["E18"] = '  -0.42%  '
$updates["D19"] = '243.38'
$updates["E19"] = '  -0.78%  '
$updates["E20"] = '  -2.44%  '
$updates["D21"] = '11.24'
$updates["E21"] = '  -2.36%  '
$updates["E22"] = '  -0.27%  '
$updates["D23"] = '4.11'
$updates["E23"] = '  -1.32%  '
$updates["E24"] = '  +4.35%  '
$updates["D25"] = '171.69'
$updates["E25"] = '  -1.36%  '
$updates["D26"] = '7.84'
$updates["E26"] = '  +0.76%  '
$updates["D27"] = '17.28'
$updates["E27"] = '  +2.86%  '
$updates["D29"] = '1.00'
$updates["E29"] = '  -0.25%  '
$updates["D30"] = '3.83'
$updates["E30"] = '  +0.13%  '
$updates["E31"] = '  -2.17%  '
$updates["E32"] = '  -1.07%  '
$updates["E33"] = '  -2.54%  '
$updates["E34"] = '  -1.14%  '
$updates["D35"] = '1.365.69'
$updates["E35"] = '  -2.10%  '
$updates["E36"] = '  -4.04%  '
$updates["E37"] = '  +0.25%  '
$updates["D38"] = '2.38'
$updates["E38"] = '  -5.15%  '
$updates["E39"] = '  -1.41%  '
$updates["D40"] = '2.42'
$updates["E40"] = '  +0.58%  '
$updates["D41"] = '81.51'
$updates["E41"] = '  -2.39%  '
$updates["E42"] = '  -1.65%  '
$updates["D43"] = '0.941'
$updates["E43"] = '  -0.83%  '
$updates["B44"] = 'InjectiveProtocol'
$updates["C44"] = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$updates["D44"] = '13.73'
$updates["E44"] = '  +0.84%  '
$updates["B45"] = 'WEMIXToken'
$updates["C45"] = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$updates["D45"] = '1.16'
$updates["E45"] = '  +3.95%  '
$updates["D46"] = '0.0503'
$updates["E46"] = '  -1.42%  '
$updates["D47"] = '1.973.67'
$updates["E47"] = '  +0.55%  '
$updates["E48"] = '  -2.31%  '
$updates["D49"] = '1.00'
$updates["E49"] = '  -0.25%  '
$updates["D50"] = '102.99'
$updates["E50"] = '  -1.93%  '
$updates["D51"] = '0.0₆0120'
$updates["E51"] = '  -7.59%  '

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

